$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '23.609.51'
$ws.Range("E2").Value = '  +0.76%  '
$ws.Range("D3").Value = '1.646.31'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").Value = '0.9979'
$ws.Range("E5").Value = '  -0.36%  '
$ws.Range("D6").Value = '304.00'
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").Value = '0.3800'
$ws.Range("E7").Value = '  +0.50%  '
$ws.Range("D8").Value = '51.93'
$ws.Range("E8").Value = '  +0.96%  '
$ws.Range("D9").Value = '0.3604'
$ws.Range("E9").Value = '  -1.37%  '
$ws.Range("D10").Value = '1.243'
$ws.Range("E10").Value = '  +0.95%  '
$ws.Range("D11").Value = '0.08192'
$ws.Range("E11").Value = '  -0.69%  '
$ws.Range("D12").Value = '0.9990'
$ws.Range("E12").Value = '  -0.25%  '
$ws.Range("D13").Value = '22.46'
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").Value = '6.518'
$ws.Range("E14").Value = '  -0.68%  '
$ws.Range("D15").Value = '7.367'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("D16").Value = '0.00001229'
$ws.Range("E16").Value = '  -1.98%  '
$ws.Range("D17").Value = '1.649.84'
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").Value = '96.90'
$ws.Range("E18").Value = '  +2.93%  '
$ws.Range("D19").Value = '0.06978'
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '6.731'
$ws.Range("E20").Value = '  +3.00%  '
$ws.Range("D21").Value = '17.57'
$ws.Range("E21").Value = '  -0.88%  '
$ws.Range("D22").Value = '0.9974'
$ws.Range("E22").Value = '  -0.43%  '
$ws.Range("D23").Value = '12.55'
$ws.Range("E23").Value = '  -1.52%  '
$ws.Range("D24").Value = '23.619.45'
$ws.Range("E24").Value = '  +0.79%  '
$ws.Range("D25").Value = '2.523'
$ws.Range("E25").Value = '  +2.12%  '
$ws.Range("D26").Value = '3.117'
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").Value = '21.27'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = '152.29'
$ws.Range("E28").Value = '  +1.11%  '
$ws.Range("D29").Value = '5.202'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").Value = '134.93'
$ws.Range("E30").Value = '  +0.37%  '
$ws.Range("D31").Value = '1.829.49'
$ws.Range("E31").Value = '  +1.10%  '
$ws.Range("D32").Value = '6.759'
$ws.Range("E32").Value = '  -1.13%  '
$ws.Range("D33").Value = '1.088'
$ws.Range("E33").Value = '  +6.60%  '
$ws.Range("D34").Value = '11.61'
$ws.Range("E34").Value = '  +6.37%  '
$ws.Range("D35").Value = '2.049'
$ws.Range("E35").Value = '  -9.22%  '
$ws.Range("D36").Value = '0.02801'
$ws.Range("E36").Value = '  +0.10%  '
$ws.Range("D37").Value = '0.2510'
$ws.Range("E37").Value = '  -0.79%  '
$ws.Range("D38").Value = '0.08815'
$ws.Range("E38").Value = '  +0.33%  '
$ws.Range("D39").Value = '6.090'
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '0.07024'
$ws.Range("E40").Value = '  -1.56%  '
$ws.Range("D41").Value = '12.82'
$ws.Range("E41").Value = '  +4.90%  '
$ws.Range("D42").Value = '0.7046'
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").Value = '1.330'
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").Value = '15.83'
$ws.Range("E44").Value = '  -1.56%  '
$ws.Range("D45").Value = '0.6498'
$ws.Range("E45").Value = '  -1.25%  '
$ws.Range("D46").Value = '2.333'
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("D47").Value = '0.9980'
$ws.Range("E47").Value = '  -0.27%  '
$ws.Range("D48").Value = '3.977'
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("D49").Value = '0.07979'
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").Value = '127.78'
$ws.Range("E50").Value = '  +1.02%  '
$ws.Range("D51").Value = '1.188'
$ws.Range("E51").Value = '  -0.97%  '
